# Adds the "ribbon reference" mini-tables and extra sample data that were
# appended below the existing Sheet1 content (rows 25-48), as well as a
# couple of "Data" labels next to the existing Fill-Handle rows (25-27).
#
# NOTE: the Value assignments below are intentionally ordered so that new
# shared-string entries are first introduced in the same sequence as in the
# target workbook (Excel's shared string table is append-only in first-seen
# order), even though that means cells aren't filled in strict reading
# order on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 25-27: existing "Double Click to AutoFill" rows gain a "Data" tag
# in column H (string already present elsewhere in the workbook).
$ws.Range("H25").Value = "Data"
$ws.Range("H26").Value = "Data"
$ws.Range("H27").Value = "Data"

# Row 28: new sample number plus the start of the "Fill Series" mini table.
$ws.Range("D28").Value = 748
$ws.Range("G28").Value = "Fill Series"

# Row 30: "Covert to Range" mini table (label first, ribbon path second).
$ws.Range("G30").Value = "Covert to Range"

# Row 32: "Conditional Formatting" mini table.
$ws.Range("G32").Value = "Conditional Formatting "
$ws.Range("H32").Value = "Highlight Cell Rules"
$ws.Range("I32").Value = "Duplicate Values"

# Back to row 28: ribbon path for Fill Series (Home > Fill), plus the
# reused "Series" label.
$ws.Range("H28").Value = "Home"
$ws.Range("I28").Value = "Fill"
$ws.Range("J28").Value = "Series"

# Row 30: ribbon path for Covert to Range (reuses existing "Table Design").
$ws.Range("H30").Value = "Table Design"

# Rows 34-36: text function reference list.
$ws.Range("G34").Value = "LEN"
$ws.Range("G35").Value = "CONCAT"
$ws.Range("G36").Value = "TEXT JOIN"

# Rows 38-40: date function reference list.
$ws.Range("G38").Value = "TODAY"
$ws.Range("G39").Value = "WEEKDAY"
$ws.Range("G40").Value = "WORKDAY"

# Rows 42-43: Protect Worksheet / Find & Replace reference entries.
$ws.Range("G42").Value = "Protect Worksheet"
$ws.Range("G43").Value = "Find & Replace"
$ws.Range("H42").Value = "Review"

# Row 37: TRIM function, added to the text list after the date block.
$ws.Range("G37").Value = "TRIM"

# Rows 45-48: logical function reference list.
$ws.Range("G45").Value = "AND"
$ws.Range("G46").Value = "OR"
$ws.Range("G47").Value = "NOT"
$ws.Range("G48").Value = "IF"

# Leave the selection where the author's workbook ended up.
$ws.Range("G49").Select()
